# Updated the log until 17th march
#
# - Adds an hours entry on 17 March (column W) to the existing
#   "DB: Implementation of database in java." row (row 26).
# - Inserts four new task rows (30-33) for work done on 17 March, each
#   logging 6 hours in column W:
#     KTN: Functionality for GUI class
#     Java fx: Main screen
#     Java fx: Viewscreen and logic   (set after "Mail class:" so the
#                                      shared-string table gets the same
#                                      ordering as the source edit)
#     Mail class:
#   The insert pushes the old "blank row" (row 30) down to row 34, and the
#   trailing summary rows (old 34/35) down to rows 38/39.
# - Moves the active-cell selection to W27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the four new log rows; everything from row 30 down
# (including the old blank row and the trailing summary rows) shifts down.
$ws.Rows("30:33").Insert()

# Extra hours logged against the existing "DB: Implementation..." task.
$ws.Range("W26").Value = 6

# New tasks worked on 17 March.
$ws.Range("A30").Value = "KTN: Functionality for GUI class"
$ws.Range("W30").Value = 6

$ws.Range("A31").Value = "Java fx: Main screen"
$ws.Range("W31").Value = 6

$ws.Range("A33").Value = "Mail class:"
$ws.Range("W33").Value = 6

$ws.Range("A32").Value = "Java fx: Viewscreen and logic"
$ws.Range("W32").Value = 6

# Leave the selection where the author left it.
$ws.Range("W27").Select()
